$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4997.5
$ws.Range("I64").Value = 3490
$ws.Range("K64").Value = 3490
$ws.Range("M64").Value = -3242
$ws.Range("H67").Value = 4997.5
$ws.Range("I67").Value = 3490
$ws.Range("K67").Value = 3490
$ws.Range("M67").Value = -2632
$ws.Range("H92").Value = 930.0909
$ws.Range("I92").Value = 747.3333
$ws.Range("J92").Value = 1752.5
$ws.Range("K92").Value = 747.3333
$ws.Range("L92").Value = 1752.5
$ws.Range("M92").Value = 500.6667
$ws.Range("N92").Value = -4248.5
$ws.Range("H125").Value = 2999.5
$ws.Range("I125").Value = 1000
$ws.Range("K125").Value = 9000
$ws.Range("M125").Value = -6540
$ws.Range("H137").Value = 2385.5715
$ws.Range("I137").Value = 1999.5
$ws.Range("J137").Value = 2540
$ws.Range("K137").Value = 5998.5
$ws.Range("L137").Value = 7620
$ws.Range("M137").Value = -3448.5
$ws.Range("N137").Value = -12720
$ws.Range("H138").Value = 4990.9287
$ws.Range("I138").Value = 4309.125
$ws.Range("J138").Value = 5900
$ws.Range("K138").Value = 12927.375
$ws.Range("L138").Value = 17700
$ws.Range("M138").Value = -7787.375
$ws.Range("N138").Value = -27980
$ws.Range("H141").Value = 2407.8333
$ws.Range("I141").Value = 1236.75
$ws.Range("J141").Value = 4750
$ws.Range("K141").Value = 3710.25
$ws.Range("L141").Value = 14250
$ws.Range("M141").Value = 1469.75
$ws.Range("N141").Value = -24610

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2408503.2
$ws.Range("I32").Value = 2260390.5
$ws.Range("K32").Value = 2260390.5
$ws.Range("M32").Value = -2260103.5
$ws.Range("H46").Value = 9934.5
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 9934.5
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 9934.5
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -10572.5
$ws.Range("H61").Value = 4885.6
$ws.Range("I61").Value = 4885.6
$ws.Range("K61").Value = 4885.6
$ws.Range("M61").Value = -4673.6
$ws.Range("H74").Value = 1416
$ws.Range("I74").Value = 1416
$ws.Range("K74").Value = 1416
$ws.Range("M74").Value = -542
$ws.Range("H77").Value = 1416
$ws.Range("I77").Value = 1416
$ws.Range("K77").Value = 7080
$ws.Range("M77").Value = -2712
$ws.Range("H136").Value = 4885.6
$ws.Range("I136").Value = 4885.6
$ws.Range("K136").Value = 14656.8
$ws.Range("M136").Value = -12106.8

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1769.6666
$ws.Range("J20").Value = 2554.5
$ws.Range("L20").Value = 2554.5
$ws.Range("N20").Value = -3048.5
$ws.Range("H86").Value = 2449.75
$ws.Range("I86").Value = 1933
$ws.Range("J86").Value = 4000
$ws.Range("K86").Value = 1933
$ws.Range("L86").Value = 4000
$ws.Range("M86").Value = -810
$ws.Range("N86").Value = -6246
$ws.Range("H89").Value = 2449.75
$ws.Range("I89").Value = 1933
$ws.Range("J89").Value = 4000
$ws.Range("K89").Value = 9665
$ws.Range("L89").Value = 20000
$ws.Range("M89").Value = -4049
$ws.Range("N89").Value = -31232
$ws.Range("H99").Value = 3538.75
$ws.Range("I99").Value = 3718.3333
$ws.Range("K99").Value = 3718.3333
$ws.Range("M99").Value = -2220.3333

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2410.0605
$ws.Range("J31").Value = 2708.077
$ws.Range("L31").Value = 2708.077
$ws.Range("N31").Value = -3298.077
$ws.Range("H32").Value = 24500
$ws.Range("I32").Value = 9000
$ws.Range("K32").Value = 9000
$ws.Range("M32").Value = -8684
$ws.Range("H34").Value = 2410.0605
$ws.Range("J34").Value = 2708.077
$ws.Range("L34").Value = 2708.077
$ws.Range("N34").Value = -3112.077

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 4850
$ws.Range("J68").Value = 4850
$ws.Range("L68").Value = 14550
$ws.Range("N68").Value = -16172
$ws.Range("H71").Value = 4850
$ws.Range("J71").Value = 4850
$ws.Range("L71").Value = 43650
$ws.Range("N71").Value = -51762
$ws.Range("H107").Value = 4994.1665
$ws.Range("I107").Value = 4994
$ws.Range("J107").Value = 4994.5
$ws.Range("K107").Value = 14982
$ws.Range("L107").Value = 14983.5
$ws.Range("M107").Value = -13062
$ws.Range("N107").Value = -18823.5
$ws.Range("H113").Value = 1999.4
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1999.4
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 5998.200000000001
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -10338.2

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 7544499.5
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H80").Value = 22399.4
$ws.Range("I80").Value = 3250
$ws.Range("J80").Value = 35165.668
$ws.Range("K80").Value = 3250
$ws.Range("L80").Value = 35165.668
$ws.Range("M80").Value = -2252
$ws.Range("N80").Value = -37161.668
$ws.Range("H83").Value = 22399.4
$ws.Range("I83").Value = 3250
$ws.Range("J83").Value = 35165.668
$ws.Range("K83").Value = 16250
$ws.Range("L83").Value = 175828.34
$ws.Range("M83").Value = -11258
$ws.Range("N83").Value = -185812.34

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 90077
$ws.Range("I63").Value = 90077
$ws.Range("K63").Value = 90077
$ws.Range("M63").Value = -89328
$ws.Range("H66").Value = 90077
$ws.Range("I66").Value = 90077
$ws.Range("K66").Value = 270231
$ws.Range("M66").Value = -266487
$ws.Range("H100").Value = 1259.2858
$ws.Range("I100").Value = 1259.2858
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1259.2858
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -718.2858000000001
$ws.Range("N100").ClearContents()
$ws.Range("H132").Value = 3864
$ws.Range("I132").Value = 3915
$ws.Range("J132").Value = 3779
$ws.Range("K132").Value = 11745
$ws.Range("L132").Value = 11337
$ws.Range("M132").Value = -9215
$ws.Range("N132").Value = -16397

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3015.1538
$ws.Range("I136").Value = 2808.0833
$ws.Range("K136").Value = 8424.249899999999
$ws.Range("M136").Value = -5874.249899999999
